$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block occupies rows 2..27 (columns B..F; A is constant across all rows).
# The edit rotates the block: the last two rows (26,27 - "provisional") move to the
# top (2,3), and the rest (2..25 - "accredited") shift down by two rows (4..27).
# Capture all existing values first so writes don't clobber reads.

$firstRow = 2
$lastRow = 27

$bVals = @{}
$cVals = @{}
$dVals = @{}
$eVals = @{}
$fVals = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $bVals[$r] = $ws.Cells.Item($r, 2).Value()
    $cVals[$r] = $ws.Cells.Item($r, 3).Value()
    $dVals[$r] = $ws.Cells.Item($r, 4).Value()
    $eVals[$r] = $ws.Cells.Item($r, 5).Value()
    $fVals[$r] = $ws.Cells.Item($r, 6).Value()
}

# Build the rotated order: source row for each destination row.
$srcForDest = @{}
$srcForDest[2] = 26
$srcForDest[3] = 27
for ($destRow = 4; $destRow -le 27; $destRow++) {
    $srcForDest[$destRow] = $destRow - 2
}

for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $src = $srcForDest[$destRow]
    $ws.Cells.Item($destRow, 2).Value = $bVals[$src]
    $ws.Cells.Item($destRow, 3).Value = $cVals[$src]
    $ws.Cells.Item($destRow, 4).Value = $dVals[$src]
    $ws.Cells.Item($destRow, 5).Value = $eVals[$src]
    $ws.Cells.Item($destRow, 6).Value = $fVals[$src]
}

# Two additional single-cell corrections beyond the rotation: the
# "count_schools_double_addresses_subtotal" (column E) for the "Table 34"
# and "Table 36" rows drop from 1 to 0. After the rotation above those
# rows now live at row 5 ("Table 34") and row 3 ("Table 36").
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(3, 5).Value = 0

# Update the last-selected cell to match the saved session state.
$ws.Range("B8").Select()
